$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accidentes_MEX")

# Update data values
$ws.Range("B2").Value = 331938
$ws.Range("B20").Value = 355281
$ws.Range("B22").Value = 315068

# Update the view: zoom level and selection range
$ws.Activate()
$excel.ActiveWindow.Zoom = 126
$ws.Range("A2:B20").Select()
